# validacao matriz na criacao de usuarios
#
# - "profiles": add three new profile rows (sava/tecnico, teste/teste,
#   novo_sistema/novo_perfil), all with description "teste" (used to
#   exercise the new user-creation validation).
# - "matriz" (the Segregation-of-Duties conflict matrix): trim back down
#   to just the single sava/aluno <-> sava/professor conflict pair.
# - "users": change the existing aluno/professor pair for cpf 1 so cpf 1
#   is a "professor", then add a new cpf 2 with both "aluno" and
#   "professor" roles on "sava" (the conflicting pair from the matriz),
#   to validate the SoD check at user-creation time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# profiles: append rows 6-8
# ---------------------------------------------------------------------
$profiles = $wb.Worksheets.Item("profiles")

$profiles.Range("A6").Value = "sava"
$profiles.Range("B6").Value = "tecnico"
$profiles.Range("C6").Value = "teste"

$profiles.Range("A7").Value = "teste"
$profiles.Range("B7").Value = "teste"
$profiles.Range("C7").Value = "teste"

$profiles.Range("A8").Value = "novo_sistema"
$profiles.Range("B8").Value = "novo_perfil"
$profiles.Range("C8").Value = "teste"

# ---------------------------------------------------------------------
# matriz: drop rows 3 and 4, keeping only the header + the sava/aluno
# vs sava/professor conflict
# ---------------------------------------------------------------------
$matriz = $wb.Worksheets.Item("matriz")
$matriz.Rows("3:4").Delete()

# ---------------------------------------------------------------------
# users: cpf 1 becomes "professor"; insert a new cpf 2 "aluno" row, and
# re-point the former row 3 (cpf 1 / professor, stored as text "1") to
# cpf 2 (as text "2") on row 4
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("users")

$users.Range("C2").Value = "professor"

$users.Rows("3:3").Insert()
$users.Range("A3").Value = 2
$users.Range("B3").Value = "sava"
$users.Range("C3").Value = "aluno"

$users.Range("A4").NumberFormat = "@"
$users.Range("A4").Value = "2"
